$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 62, shifting existing rows 62-171 down to 63-172.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data record.
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 45100
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100108
$ws.Cells.Item(62, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value = 100108002
$ws.Cells.Item(62, 10).Value = "Mango"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 80
$ws.Cells.Item(62, 14).Value = 9000
$ws.Cells.Item(62, 15).Value = 9000
$ws.Cells.Item(62, 16).Value = 9000
$ws.Cells.Item(62, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(62, 18).Value = "Ecuador"
$ws.Cells.Item(62, 19).Value = 2250
$ws.Cells.Item(62, 20).Value = 4
